$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> [new Price (D) text, new Volume(1h) (E) text]; empty string means "unchanged"
$updates = @{
    2 = @('23.888.19', '  -2.09%  ')
    3 = @('1.651.31', '  -1.01%  ')
    4 = @('1.001', '  -0.04%  ')
    5 = @('310.99', '  -0.38%  ')
    6 = @('', '  -0.04%  ')
    7 = @('0.3895', '  -1.28%  ')
    8 = @('0.3823', '  -2.43%  ')
    9 = @('51.94', '  -0.41%  ')
    10 = @('1.342', '  -3.25%  ')
    11 = @('1.001', '  -0.05%  ')
    12 = @('0.08449', '  -1.40%  ')
    13 = @('23.91', '  -2.10%  ')
    14 = @('7.014', '  -3.63%  ')
    15 = @('8.021', '  +0.33%  ')
    16 = @('0.00001314', '  -1.44%  ')
    17 = @('1.649.81', '  -0.87%  ')
    18 = @('94.13', '  -0.97%  ')
    19 = @('0.06984', '  -0.71%  ')
    20 = @('19.57', '  -4.68%  ')
    21 = @('6.957', '  -0.43%  ')
    22 = @('', '  +0.03%  ')
    23 = @('13.67', '  -0.51%  ')
    24 = @('23.890.44', '  -2.08%  ')
    25 = @('2.445', '  -2.99%  ')
    26 = @('2.935', '  -4.82%  ')
    27 = @('22.03', '  -2.25%  ')
    28 = @('153.44', '  -2.22%  ')
    29 = @('5.399', '  +0.93%  ')
    30 = @('137.29', '  -3.68%  ')
    31 = @('7.742', '  -2.26%  ')
    32 = @('2.484', '  -2.31%  ')
    33 = @('1.833.23', '')
    34 = @('0.08129', '  -1.15%  ')
    35 = @('0.9987', '  -5.63%  ')
    36 = @('6.666', '  -3.02%  ')
    37 = @('', '  -6.96%  ')
    38 = @('', '  -3.18%  ')
    39 = @('10.71', '  -4.44%  ')
    40 = @('0.09102', '  -1.86%  ')
    41 = @('0.7553', '  -1.68%  ')
    42 = @('', '  -1.90%  ')
    43 = @('1.418', '  -1.92%  ')
    44 = @('16.45', '  -0.33%  ')
    45 = @('0.6940', '  -2.04%  ')
    46 = @('', '  -3.61%  ')
    47 = @('4.108', '  -0.47%  ')
    48 = @('1.001', '  -0.05%  ')
    49 = @('0.08272', '  -1.78%  ')
    50 = @('133.51', '')
    51 = @('1.226', '  -2.74%  ')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dNew = $vals[0]
    $eNew = $vals[1]
    if ($dNew -ne "") {
        # Force the cell to remain plain text (matches original inline-string cells)
        # even when the new value looks like a number (e.g. "1.001").
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $dNew
        $cell.Style = "Normal"
    }
    if ($eNew -ne "") {
        $ws.Range("E$row").Value = $eNew
    }
}
